# Update Il4-Il4ra NATMI edge-weight metrics following Dr Hou advice
# (recomputed ligand/receptor expressing-cell counts and derived stats)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7922663333333334
$ws.Range("H2").Value = 2.376799
$ws.Range("I2").Value = 0.1759587713796512
$ws.Range("J2").Value = 0.1759587713796512
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.41769233333333
$ws.Range("N2").Value = 106.253077
$ws.Range("O2").Value = 0.3584901985342537
$ws.Range("P2").Value = 0.3584901985342537
$ws.Range("Q2").Value = 28.06024524005811
$ws.Range("R2").Value = 252.542207160523
$ws.Range("S2").Value = 0.06307949488573451
$ws.Range("T2").Value = 0.06307949488573453

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7922663333333334
$ws.Range("H3").Value = 2.376799
$ws.Range("I3").Value = 0.1759587713796512
$ws.Range("J3").Value = 0.1759587713796512
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.255923
$ws.Range("N3").Value = 90.767769
$ws.Range("O3").Value = 0.3062438890999955
$ws.Range("P3").Value = 0.3062438890999956
$ws.Range("Q3").Value = 23.97074917682567
$ws.Range("R3").Value = 215.736742591431
$ws.Range("S3").Value = 0.05388629846856136
$ws.Range("T3").Value = 0.05388629846856138

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.7922663333333334
$ws.Range("H4").Value = 2.376799
$ws.Range("I4").Value = 0.1759587713796512
$ws.Range("J4").Value = 0.1759587713796512
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 29.46429466666667
$ws.Range("N4").Value = 88.39288400000001
$ws.Range("O4").Value = 0.2982311988402488
$ws.Range("P4").Value = 0.2982311988402488
$ws.Range("Q4").Value = 23.3435686998129
$ws.Range("R4").Value = 210.092118298316
$ws.Range("S4").Value = 0.05247639533501064
$ws.Range("T4").Value = 0.05247639533501065

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7922663333333334
$ws.Range("H5").Value = 2.376799
$ws.Range("I5").Value = 0.1759587713796512
$ws.Range("J5").Value = 0.1759587713796512
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.658912
$ws.Range("N5").Value = 10.976736
$ws.Range("O5").Value = 0.03703471352550186
$ws.Range("P5").Value = 0.03703471352550187
$ws.Range("Q5").Value = 2.898832794229334
$ws.Range("R5").Value = 26.089495148064
$ws.Range("S5").Value = 0.006516582690344657
$ws.Range("T5").Value = 0.00651658269034466

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.077831666666667
$ws.Range("H6").Value = 6.233495
$ws.Range("I6").Value = 0.4614770208171574
$ws.Range("J6").Value = 0.4614770208171574
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.41769233333333
$ws.Range("N6").Value = 106.253077
$ws.Range("O6").Value = 0.3584901985342537
$ws.Range("P6").Value = 0.3584901985342537
$ws.Range("Q6").Value = 73.59200269045722
$ws.Range("R6").Value = 662.328024214115
$ws.Range("S6").Value = 0.1654349888117387
$ws.Range("T6").Value = 0.1654349888117387

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.077831666666667
$ws.Range("H7").Value = 6.233495
$ws.Range("I7").Value = 0.4614770208171574
$ws.Range("J7").Value = 0.4614770208171574
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 30.255923
$ws.Range("N7").Value = 90.767769
$ws.Range("O7").Value = 0.3062438890999955
$ws.Range("P7").Value = 0.3062438890999956
$ws.Range("Q7").Value = 62.86671491362834
$ws.Range("R7").Value = 565.8004342226551
$ws.Range("S7").Value = 0.1413245175853259
$ws.Range("T7").Value = 0.1413245175853259

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.077831666666667
$ws.Range("H8").Value = 6.233495
$ws.Range("I8").Value = 0.4614770208171574
$ws.Range("J8").Value = 0.4614770208171574
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.46429466666667
$ws.Range("N8").Value = 88.39288400000001
$ws.Range("O8").Value = 0.2982311988402488
$ws.Range("P8").Value = 0.2982311988402488
$ws.Range("Q8").Value = 61.22184449439779
$ws.Range("R8").Value = 550.9966004495801
$ws.Range("S8").Value = 0.1376268451555273
$ws.Range("T8").Value = 0.1376268451555273

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.077831666666667
$ws.Range("H9").Value = 6.233495
$ws.Range("I9").Value = 0.4614770208171574
$ws.Range("J9").Value = 0.4614770208171574
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.658912
$ws.Range("N9").Value = 10.976736
$ws.Range("O9").Value = 0.03703471352550186
$ws.Range("P9").Value = 0.03703471352550187
$ws.Range("Q9").Value = 7.602603219146668
$ws.Range("R9").Value = 68.42342897232001
$ws.Range("S9").Value = 0.01709066926456548
$ws.Range("T9").Value = 0.01709066926456549

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8480786666666668
$ws.Range("H10").Value = 2.544236
$ws.Range("I10").Value = 0.1883544383264543
$ws.Range("J10").Value = 0.1883544383264543
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 35.41769233333333
$ws.Range("N10").Value = 106.253077
$ws.Range("O10").Value = 0.3584901985342537
$ws.Range("P10").Value = 0.3584901985342537
$ws.Range("Q10").Value = 30.03698929046356
$ws.Range("R10").Value = 270.332903614172
$ws.Range("S10").Value = 0.06752321999045845
$ws.Range("T10").Value = 0.06752321999045845

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8480786666666668
$ws.Range("H11").Value = 2.544236
$ws.Range("I11").Value = 0.1883544383264543
$ws.Range("J11").Value = 0.1883544383264543
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 30.255923
$ws.Range("N11").Value = 90.767769
$ws.Range("O11").Value = 0.3062438890999955
$ws.Range("P11").Value = 0.3062438890999956
$ws.Range("Q11").Value = 25.65940283660933
$ws.Range("R11").Value = 230.934625529484
$ws.Range("S11").Value = 0.0576823957223386
$ws.Range("T11").Value = 0.05768239572233862

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.8480786666666668
$ws.Range("H12").Value = 2.544236
$ws.Range("I12").Value = 0.1883544383264543
$ws.Range("J12").Value = 0.1883544383264543
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 29.46429466666667
$ws.Range("N12").Value = 88.39288400000001
$ws.Range("O12").Value = 0.2982311988402488
$ws.Range("P12").Value = 0.2982311988402488
$ws.Range("Q12").Value = 24.98803973518045
$ws.Range("R12").Value = 224.892357616624
$ws.Range("S12").Value = 0.05617316994898017
$ws.Range("T12").Value = 0.05617316994898018

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.8480786666666668
$ws.Range("H13").Value = 2.544236
$ws.Range("I13").Value = 0.1883544383264543
$ws.Range("J13").Value = 0.1883544383264543
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.658912
$ws.Range("N13").Value = 10.976736
$ws.Range("O13").Value = 0.03703471352550186
$ws.Range("P13").Value = 0.03703471352550187
$ws.Range("Q13").Value = 3.103045210410667
$ws.Range("R13").Value = 27.927406893696
$ws.Range("S13").Value = 0.006975652664677042
$ws.Range("T13").Value = 0.006975652664677044

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.7843913333333336
$ws.Range("H14").Value = 2.353174000000001
$ws.Range("I14").Value = 0.1742097694767371
$ws.Range("J14").Value = 0.1742097694767372
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 35.41769233333333
$ws.Range("N14").Value = 106.253077
$ws.Range("O14").Value = 0.3584901985342537
$ws.Range("P14").Value = 0.3584901985342537
$ws.Range("Q14").Value = 27.78133091293311
$ws.Range("R14").Value = 250.031978216398
$ws.Range("S14").Value = 0.06245249484632207
$ws.Range("T14").Value = 0.06245249484632208

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.7843913333333336
$ws.Range("H15").Value = 2.353174000000001
$ws.Range("I15").Value = 0.1742097694767371
$ws.Range("J15").Value = 0.1742097694767372
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 30.255923
$ws.Range("N15").Value = 90.767769
$ws.Range("O15").Value = 0.3062438890999955
$ws.Range("P15").Value = 0.3062438890999956
$ws.Range("Q15").Value = 23.73248378320067
$ws.Range("R15").Value = 213.592354048806
$ws.Range("S15").Value = 0.05335067732376967
$ws.Range("T15").Value = 0.05335067732376969

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.7843913333333336
$ws.Range("H16").Value = 2.353174000000001
$ws.Range("I16").Value = 0.1742097694767371
$ws.Range("J16").Value = 0.1742097694767372
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 29.46429466666667
$ws.Range("N16").Value = 88.39288400000001
$ws.Range("O16").Value = 0.2982311988402488
$ws.Range("P16").Value = 0.2982311988402488
$ws.Range("Q16").Value = 23.1115373793129
$ws.Range("R16").Value = 208.0038364138161
$ws.Range("S16").Value = 0.05195478840073071
$ws.Range("T16").Value = 0.05195478840073071

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.7843913333333336
$ws.Range("H17").Value = 2.353174000000001
$ws.Range("I17").Value = 0.1742097694767371
$ws.Range("J17").Value = 0.1742097694767372
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.658912
$ws.Range("N17").Value = 10.976736
$ws.Range("O17").Value = 0.03703471352550186
$ws.Range("P17").Value = 0.03703471352550187
$ws.Range("Q17").Value = 2.870018862229335
$ws.Range("R17").Value = 25.83016976006401
$ws.Range("S17").Value = 0.006451808905914678
$ws.Range("T17").Value = 0.00645180890591468
